$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,7).Value = 74.609651
$ws.Cells.Item(2,8).Value = 223.828953
$ws.Cells.Item(2,9).Value = 0.1061386348809139
$ws.Cells.Item(2,10).Value = 0.1061386348809139
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 0.6763629999999999
$ws.Cells.Item(2,14).Value = 2.029089
$ws.Cells.Item(2,15).Value = 0.6127318215515719
$ws.Cells.Item(2,16).Value = 0.6127318215515719
$ws.Cells.Item(2,17).Value = 50.463207379313
$ws.Cells.Item(2,18).Value = 454.168866413817
$ws.Cells.Item(2,19).Value = 0.06503451908757958
$ws.Cells.Item(2,20).Value = 0.06503451908757961
$ws.Cells.Item(3,7).Value = 74.609651
$ws.Cells.Item(3,8).Value = 223.828953
$ws.Cells.Item(3,9).Value = 0.1061386348809139
$ws.Cells.Item(3,10).Value = 0.1061386348809139
$ws.Cells.Item(3,13).Value = 0.05377866666666667
$ws.Cells.Item(3,15).Value = 0.04871925339984812
$ws.Cells.Item(3,16).Value = 0.04871925339984811
$ws.Cells.Item(3,17).Value = 4.012407551245333
$ws.Cells.Item(3,18).Value = 36.11166796120801
$ws.Cells.Item(3,19).Value = 0.005170995048277203
$ws.Cells.Item(3,20).Value = 0.005170995048277204
$ws.Cells.Item(4,7).Value = 74.609651
$ws.Cells.Item(4,8).Value = 223.828953
$ws.Cells.Item(4,9).Value = 0.1061386348809139
$ws.Cells.Item(4,10).Value = 0.1061386348809139
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 0.3737066666666666
$ws.Cells.Item(4,14).Value = 1.12112
$ws.Cells.Item(4,15).Value = 0.3385489250485801
$ws.Cells.Item(4,16).Value = 0.33854892504858
$ws.Cells.Item(4,17).Value = 27.88212397637333
$ws.Cells.Item(4,18).Value = 250.93911578736
$ws.Cells.Item(4,19).Value = 0.03593312074505713
$ws.Cells.Item(4,20).Value = 0.03593312074505713
$ws.Cells.Item(5,7).Value = 597.374756
$ws.Cells.Item(5,9).Value = 0.8498168837991085
$ws.Cells.Item(5,10).Value = 0.8498168837991086
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.6763629999999999
$ws.Cells.Item(5,14).Value = 2.029089
$ws.Cells.Item(5,15).Value = 0.6127318215515719
$ws.Cells.Item(5,16).Value = 0.6127318215515719
$ws.Cells.Item(5,17).Value = 404.042182092428
$ws.Cells.Item(5,18).Value = 3636.379638831852
$ws.Cells.Item(5,19).Value = 0.5207098471955083
$ws.Cells.Item(5,20).Value = 0.5207098471955083
$ws.Cells.Item(6,7).Value = 597.374756
$ws.Cells.Item(6,9).Value = 0.8498168837991085
$ws.Cells.Item(6,10).Value = 0.8498168837991086
$ws.Cells.Item(6,13).Value = 0.05377866666666667
$ws.Cells.Item(6,15).Value = 0.04871925339984812
$ws.Cells.Item(6,16).Value = 0.04871925339984811
$ws.Cells.Item(6,17).Value = 32.12601787800534
$ws.Cells.Item(6,18).Value = 289.134160902048
$ws.Cells.Item(6,19).Value = 0.04140244410527805
$ws.Cells.Item(6,20).Value = 0.04140244410527805
$ws.Cells.Item(7,7).Value = 597.374756
$ws.Cells.Item(7,9).Value = 0.8498168837991085
$ws.Cells.Item(7,10).Value = 0.8498168837991086
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 0.3737066666666666
$ws.Cells.Item(7,14).Value = 1.12112
$ws.Cells.Item(7,15).Value = 0.3385489250485801
$ws.Cells.Item(7,16).Value = 0.33854892504858
$ws.Cells.Item(7,17).Value = 223.2429288155733
$ws.Cells.Item(7,18).Value = 2009.18635934016
$ws.Cells.Item(7,19).Value = 0.2877045924983223
$ws.Cells.Item(7,20).Value = 0.2877045924983223
$ws.Cells.Item(8,7).Value = 30.48438
$ws.Cells.Item(8,8).Value = 91.45313999999999
$ws.Cells.Item(8,9).Value = 0.04336664808137267
$ws.Cells.Item(8,10).Value = 0.04336664808137267
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 0.6763629999999999
$ws.Cells.Item(8,14).Value = 2.029089
$ws.Cells.Item(8,15).Value = 0.6127318215515719
$ws.Cells.Item(8,16).Value = 0.6127318215515719
$ws.Cells.Item(8,17).Value = 20.61850670994
$ws.Cells.Item(8,18).Value = 185.56656038946
$ws.Cells.Item(8,19).Value = 0.02657212527348546
$ws.Cells.Item(8,20).Value = 0.02657212527348546
$ws.Cells.Item(9,7).Value = 30.48438
$ws.Cells.Item(9,8).Value = 91.45313999999999
$ws.Cells.Item(9,9).Value = 0.04336664808137267
$ws.Cells.Item(9,10).Value = 0.04336664808137267
$ws.Cells.Item(9,13).Value = 0.05377866666666667
$ws.Cells.Item(9,15).Value = 0.04871925339984812
$ws.Cells.Item(9,16).Value = 0.04871925339984811
$ws.Cells.Item(9,17).Value = 1.63940931056
$ws.Cells.Item(9,18).Value = 14.75468379504
$ws.Cells.Item(9,19).Value = 0.002112790716978432
$ws.Cells.Item(9,20).Value = 0.002112790716978432
$ws.Cells.Item(10,7).Value = 30.48438
$ws.Cells.Item(10,8).Value = 91.45313999999999
$ws.Cells.Item(10,9).Value = 0.04336664808137267
$ws.Cells.Item(10,10).Value = 0.04336664808137267
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 0.3737066666666666
$ws.Cells.Item(10,14).Value = 1.12112
$ws.Cells.Item(10,15).Value = 0.3385489250485801
$ws.Cells.Item(10,16).Value = 0.33854892504858
$ws.Cells.Item(10,17).Value = 11.3922160352
$ws.Cells.Item(10,18).Value = 102.5299443168
$ws.Cells.Item(10,19).Value = 0.01468173209090878
$ws.Cells.Item(10,20).Value = 0.01468173209090878
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 0.4764796666666666
$ws.Cells.Item(11,8).Value = 1.429439
$ws.Cells.Item(11,9).Value = 0.0006778332386049212
$ws.Cells.Item(11,10).Value = 0.0006778332386049213
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 0.6763629999999999
$ws.Cells.Item(11,14).Value = 2.029089
$ws.Cells.Item(11,15).Value = 0.6127318215515719
$ws.Cells.Item(11,16).Value = 0.6127318215515719
$ws.Cells.Item(11,17).Value = 0.3222732167856666
$ws.Cells.Item(11,18).Value = 2.900458951071
$ws.Cells.Item(11,19).Value = 0.0004153299949985946
$ws.Cells.Item(11,20).Value = 0.0004153299949985947
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 0.4764796666666666
$ws.Cells.Item(12,8).Value = 1.429439
$ws.Cells.Item(12,9).Value = 0.0006778332386049212
$ws.Cells.Item(12,10).Value = 0.0006778332386049213
$ws.Cells.Item(12,13).Value = 0.05377866666666667
$ws.Cells.Item(12,15).Value = 0.04871925339984812
$ws.Cells.Item(12,16).Value = 0.04871925339984811
$ws.Cells.Item(12,17).Value = 0.02562444116711111
$ws.Cells.Item(12,18).Value = 0.230619970504
$ws.Cells.Item(12,19).Value = 0.00003302352931443287
$ws.Cells.Item(12,20).Value = 0.00003302352931443287
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 0.4764796666666666
$ws.Cells.Item(13,8).Value = 1.429439
$ws.Cells.Item(13,9).Value = 0.0006778332386049212
$ws.Cells.Item(13,10).Value = 0.0006778332386049213
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 0.3737066666666666
$ws.Cells.Item(13,14).Value = 1.12112
$ws.Cells.Item(13,15).Value = 0.3385489250485801
$ws.Cells.Item(13,16).Value = 0.33854892504858
$ws.Cells.Item(13,17).Value = 0.1780636279644444
$ws.Cells.Item(13,18).Value = 1.60257265168
$ws.Cells.Item(13,19).Value = 0.0002294797142918938
$ws.Cells.Item(13,20).Value = 0.0002294797142918937
